$d = $word.ActiveDocument

# Locate the "GITHUB PAGE URL : [  ]" paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*GITHUB*URL*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'GITHUB PAGE URL' paragraph"
}

# The two following (empty) paragraphs get folded into this edit too:
#  - the first is dropped entirely (just a bare new paragraph mark remains)
#  - the second keeps its right-alignment but loses its language run property
$after1 = $target.Next()
$after2 = $after1.Next()

$editRange = $d.Range($target.Range.Start, $after2.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:r><w:t xml:space="preserve">GITHUB </w:t></w:r>' +
'<w:r><w:t>PAGE</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:proofErr w:type="gramStart"/>' +
'<w:r><w:t>URL</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t>:</w:t></w:r>' +
'<w:proofErr w:type="gramEnd"/>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t>https://misadelgado10.github.io/DesarrolloWeb/Laboratorio%201/page/index2.html</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'</w:p>' +
'<w:p/>' +
'<w:p><w:pPr><w:jc w:val="right"/></w:pPr></w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$editRange.InsertXML($xml)
